# Fruta / hortaliza, semanal
# A new weekly price record (2022-03-03) is inserted for
# "Comercializadora del Agro de Limarí - Tuna", pushing the existing
# data rows down by two rows. The new record reuses the Especial /
# Primera rows that were previously at row 21-22 (2022-03-02) as a
# template, but with an updated date and volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing rows 21:22 (Especial / Primera for 2022-03-02)
# and insert the copies above themselves. This shifts the old rows
# (and everything below) down by two rows, ending up at rows 23:24.
$ws.Range("A21:T22").Copy()
$ws.Range("A21:T22").Insert()

# Update the newly inserted rows 21:22 with the new date and volumes
# for the new weekly entry (2022-03-03).
$newDate = Get-Date -Year 2022 -Month 3 -Day 3 -Hour 0 -Minute 0 -Second 0

$ws.Range("D21").Value = $newDate
$ws.Range("M21").Value = 400

$ws.Range("D22").Value = $newDate
$ws.Range("M22").Value = 400
